# 01.04.21 debugging (copying peak table)
#
# The occupancy/fragmentation peak table on sheet "analysis" was recomputed
# and re-pasted from the analysis tool. This refreshes the run timestamp in
# A1 and rewrites the "c" (column C) / "y" (column D) occupancy values for
# rows 14-27. As part of the re-paste, C19 and D15 ended up without values
# (those cells are fully cleared, not just zeroed), matching a table that
# was copied in shifted by one row for a couple of cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp in A1
$ws.Range("A1").Value2 = "01/04/2021 11:01"

# Summary figures near the top of the sheet
$ws.Range("C3").Value2 = 0.4615051419869385
$ws.Range("C4").Value2 = 0.5384948580130613

# Peak table, column C ("c" fragment occupancy) and column D ("y" fragment occupancy)
$ws.Range("C14").Value2 = 0.02000087177686502

$ws.Range("C15").Value2 = 0.3351697413725013
$ws.Range("D15").Clear()

$ws.Range("C16").Value2 = 0.4721014773042254
$ws.Range("D16").Value2 = 0.4600179281589262

$ws.Range("C17").Value2 = 0.3554065646873324
$ws.Range("D17").Value2 = 0.5027834654089542

$ws.Range("C18").Value2 = 0.4234487810123325
$ws.Range("D18").Value2 = 0.5317488919926445

$ws.Range("C19").Clear()
$ws.Range("D19").Value2 = 0.4371745787899363

$ws.Range("C20").Value2 = 0.4358916870624788
$ws.Range("D20").Value2 = 0.4689236312318184

$ws.Range("C21").Value2 = 0.5452501465705899
$ws.Range("D21").Value2 = 0.3813157300431642

$ws.Range("C22").Value2 = 0.7626031934592458
$ws.Range("D22").Value2 = 0.2165900832457863

$ws.Range("C23").Value2 = 0.8227067731027261
$ws.Range("D23").Value2 = 0.1274002410339445

$ws.Range("C24").Value2 = 0.7303401947315076
$ws.Range("D24").Value2 = 0.2021415448718557

$ws.Range("C25").Value2 = 0.8954395298949853
$ws.Range("D25").Value2 = 0.1289371743655314

$ws.Range("C26").Value2 = 0.9253528371680027
$ws.Range("D26").Value2 = 0.0327568583004501

$ws.Range("C27").Value2 = 0.941606353845397
$ws.Range("D27").Value2 = 0.01807710150318315
